# Replace missing (0) pka_acid values with the column mean of the
# non-missing values, then update the km_cluster (E) and db_cluster (F)
# columns to reflect the resulting re-clustering.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$meanPkaAcid = 8.917755102040816

# row => @{ D = pka_acid; E = km_cluster; F = db_cluster } (only changed columns listed)
$updates = @{
    2 = @{ "E" = 0 }
    3 = @{ "E" = 0 }
    4 = @{ "E" = 0 }
    5 = @{ "D" = $meanPkaAcid; "E" = 1; "F" = 2 }
    6 = @{ "D" = $meanPkaAcid }
    7 = @{ "E" = 0 }
    8 = @{ "E" = 0 }
    9 = @{ "D" = $meanPkaAcid; "F" = 2 }
    10 = @{ "D" = $meanPkaAcid }
    11 = @{ "E" = 0 }
    12 = @{ "E" = 0 }
    13 = @{ "E" = 0 }
    14 = @{ "E" = 0 }
    15 = @{ "E" = 0 }
    16 = @{ "E" = 0; "F" = 2 }
    17 = @{ "D" = $meanPkaAcid; "F" = 2 }
    18 = @{ "D" = $meanPkaAcid; "F" = 2 }
    19 = @{ "D" = $meanPkaAcid; "F" = 2 }
    20 = @{ "D" = $meanPkaAcid; "F" = 2 }
    21 = @{ "D" = $meanPkaAcid; "F" = 2 }
    22 = @{ "E" = 0 }
    23 = @{ "E" = 0 }
    24 = @{ "E" = 0 }
    25 = @{ "E" = 0 }
    26 = @{ "E" = 0 }
    27 = @{ "E" = 0 }
    28 = @{ "D" = $meanPkaAcid; "F" = 2 }
    29 = @{ "D" = $meanPkaAcid; "F" = 2 }
    30 = @{ "E" = 0 }
    31 = @{ "D" = $meanPkaAcid; "E" = 1; "F" = 2 }
    32 = @{ "D" = $meanPkaAcid; "F" = 0 }
    33 = @{ "E" = 0; "F" = 0 }
    34 = @{ "D" = $meanPkaAcid }
    35 = @{ "E" = 0 }
    36 = @{ "D" = $meanPkaAcid }
    37 = @{ "F" = 2 }
    38 = @{ "E" = 1; "F" = 2 }
    39 = @{ "E" = 1; "F" = 2 }
    40 = @{ "D" = $meanPkaAcid; "E" = 1; "F" = 2 }
    41 = @{ "D" = $meanPkaAcid; "E" = 1; "F" = 2 }
    42 = @{ "D" = $meanPkaAcid; "E" = 1; "F" = 2 }
    43 = @{ "D" = $meanPkaAcid; "E" = 1; "F" = 2 }
    44 = @{ "E" = 0; "F" = 0 }
    45 = @{ "F" = 2 }
    46 = @{ "D" = $meanPkaAcid; "E" = 1; "F" = 2 }
    47 = @{ "D" = $meanPkaAcid; "F" = 2 }
    48 = @{ "E" = 1; "F" = 2 }
    49 = @{ "D" = $meanPkaAcid; "F" = 2 }
    50 = @{ "E" = 1; "F" = 2 }
    51 = @{ "E" = 0; "F" = 2 }
    52 = @{ "F" = 2 }
    53 = @{ "D" = $meanPkaAcid; "E" = 1; "F" = 2 }
    54 = @{ "E" = 1; "F" = 2 }
    55 = @{ "F" = 2 }
    56 = @{ "E" = 1; "F" = 2 }
    57 = @{ "D" = $meanPkaAcid; "E" = 1; "F" = 2 }
    58 = @{ "E" = 0; "F" = 0 }
    59 = @{ "E" = 1; "F" = 2 }
    60 = @{ "F" = 2 }
    61 = @{ "D" = $meanPkaAcid; "E" = 1; "F" = 2 }
    62 = @{ "D" = $meanPkaAcid; "E" = 1; "F" = 2 }
    63 = @{ "E" = 0; "F" = 0 }
    64 = @{ "E" = 1; "F" = 2 }
    65 = @{ "F" = 2 }
    66 = @{ "E" = 1; "F" = 2 }
    67 = @{ "E" = 0; "F" = 2 }
    68 = @{ "E" = 1; "F" = 2 }
    69 = @{ "D" = $meanPkaAcid; "E" = 1; "F" = 2 }
    70 = @{ "F" = 2 }
    71 = @{ "E" = 1; "F" = 2 }
    72 = @{ "F" = 2 }
    73 = @{ "D" = $meanPkaAcid; "E" = 1; "F" = 2 }
    74 = @{ "D" = $meanPkaAcid; "E" = 1; "F" = 2 }
    75 = @{ "D" = $meanPkaAcid; "E" = 1; "F" = 2 }
    76 = @{ "E" = 0; "F" = 2 }
    77 = @{ "E" = 0; "F" = 2 }
    78 = @{ "D" = $meanPkaAcid; "E" = 1; "F" = 2 }
    79 = @{ "F" = 2 }
    80 = @{ "D" = $meanPkaAcid; "E" = 1; "F" = 2 }
    81 = @{ "F" = 2 }
    82 = @{ "F" = 2 }
    83 = @{ "D" = $meanPkaAcid; "F" = 2 }
    84 = @{ "D" = $meanPkaAcid; "E" = 1; "F" = 2 }
    85 = @{ "D" = $meanPkaAcid; "F" = 0 }
    86 = @{ "D" = $meanPkaAcid; "F" = 0 }
    87 = @{ "D" = $meanPkaAcid; "F" = 2 }
    88 = @{ "D" = $meanPkaAcid; "E" = 1; "F" = 2 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $colIndex = @{ "D" = 4; "E" = 5; "F" = 6 }[$col]
        $ws.Cells.Item($row, $colIndex).Value = $updates[$row][$col]
    }
}
